{"js": "// Replace the date line and each division-problem cell's text, in document\n// order, with the new values from the commit's diff. The mapping is\n// positional (old text -> new text) rather than a global find/replace,\n// since several cells share the same original text (\"77\u00f73=\") but each one\n// maps to a different replacement.\nconst replacements = [\n  \"2023-11-11 Saturday\",\n  \"34\u00f74=\",\n  \"56\u00f73=\",\n  \"43\u00f72=\",\n  \"95\u00f75=\",\n  \"52\u00f75=\",\n  \"48\u00f75=\",\n  \"12\u00f75=\",\n  \"20\u00f73=\",\n  \"14\u00f74=\",\n  \"65\u00f74=\",\n  \"90\u00f77=\",\n  \"57\u00f77=\",\n  \"80\u00f73=\",\n  \"29\u00f75=\",\n  \"45\u00f77=\",\n  \"19\u00f72=\",\n  \"45\u00f79=\",\n  \"59\u00f77=\",\n  \"26\u00f79=\",\n  \"27\u00f77=\",\n  \"10\u00f74=\",\n  \"74\u00f72=\",\n  \"15\u00f76=\",\n  \"26\u00f76=\",\n  \"57\u00f73=\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet i = 0;\nfor (const paragraph of paragraphs.items) {\n  if (i >= replacements.length) break;\n  if (paragraph.text !== \"\") {\n    paragraph.insertText(replacements[i], \"Replace\");\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and each division-problem cell's text, in document\n# order, with the new values from the commit's diff. Several cells share the\n# same original text (e.g. \"77\u00f73=\" appears three times) but each occurrence\n# maps to a different replacement, so we walk the document once, replacing\n# exactly one occurrence at a time (wdReplaceOne) and resuming the search\n# right after the text we just replaced.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2023-11-10 Friday\", \"2023-11-11 Saturday\"),\n    @(\"60\u00f77=\", \"34\u00f74=\"),\n    @(\"40\u00f75=\", \"56\u00f73=\"),\n    @(\"68\u00f73=\", \"43\u00f72=\"),\n    @(\"74\u00f72=\", \"95\u00f75=\"),\n    @(\"25\u00f77=\", \"52\u00f75=\"),\n    @(\"80\u00f78=\", \"48\u00f75=\"),\n    @(\"38\u00f79=\", \"12\u00f75=\"),\n    @(\"77\u00f73=\", \"20\u00f73=\"),\n    @(\"71\u00f78=\", \"14\u00f74=\"),\n    @(\"77\u00f73=\", \"65\u00f74=\"),\n    @(\"96\u00f78=\", \"90\u00f77=\"),\n    @(\"52\u00f74=\", \"57\u00f77=\"),\n    @(\"27\u00f77=\", \"80\u00f73=\"),\n    @(\"77\u00f73=\", \"29\u00f75=\"),\n    @(\"54\u00f72=\", \"45\u00f77=\"),\n    @(\"26\u00f73=\", \"19\u00f72=\"),\n    @(\"85\u00f77=\", \"45\u00f79=\"),\n    @(\"17\u00f72=\", \"59\u00f77=\"),\n    @(\"23\u00f79=\", \"26\u00f79=\"),\n    @(\"25\u00f76=\", \"27\u00f77=\"),\n    @(\"84\u00f78=\", \"10\u00f74=\"),\n    @(\"39\u00f79=\", \"74\u00f72=\"),\n    @(\"19\u00f72=\", \"15\u00f76=\"),\n    @(\"60\u00f75=\", \"26\u00f76=\"),\n    @(\"46\u00f76=\", \"57\u00f73=\")\n)\n\n$searchStart = 0\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $rng = $d.Range($searchStart, $d.Content.End)\n    $found = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n    if (-not $found) {\n        throw \"Could not find '$findText' starting at $searchStart\"\n    }\n    $searchStart = $rng.End\n}\n"}
